$wb = $excel.ActiveWorkbook

# Sheet "展览" — update 想去人数 (want-to-go count) for several events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3261
$ws1.Range("F9").Value = 23
$ws1.Range("F11").Value = 217
$ws1.Range("F12").Value = 1180
$ws1.Range("F13").Value = 97

# Sheet "全部类型" — same events appear here, one row lower (extra row at 8)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3261
$ws4.Range("F10").Value = 23
$ws4.Range("F12").Value = 217
$ws4.Range("F13").Value = 1180
$ws4.Range("F14").Value = 97
